$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting already used by A2:A14 onto the new rows,
# then fill in the newly-tracked time entries (date serials in column A,
# hours in column B).
$ws.Range("A14").Copy()
$ws.Range("A15:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Value = 41561
$ws.Range("B15").Value = 2.5

$ws.Range("A16").Value = 41562
$ws.Range("B16").Value = 1

$ws.Range("A17").Value = 41563
$ws.Range("B17").Value = 5

# Recompute the totals formula so the cached value reflects the new rows.
$ws.Range("B28").Formula = "=SUM(B2:B27)"

# Move the active selection as recorded in the saved view state.
$ws.Range("A18").Select()
